$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to retain text storage so numeric-looking
# strings (e.g. "1.002") are not coerced into numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '31.022.77'
$ws.Range("E2").Value = '  +0.10%  '

$ws.Range("D3").Value = '1.961.24'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '244.37'
$ws.Range("E5").Value = '  -1.34%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = '0.4850'
$ws.Range("E7").Value = '  +0.68%  '

$ws.Range("D8").Value = '0.2947'
$ws.Range("E8").Value = '  +0.45%  '

$ws.Range("D9").Value = '0.07076'
$ws.Range("E9").Value = '  +4.35%  '

$ws.Range("D10").Value = '19.70'
$ws.Range("E10").Value = '  +2.95%  '

$ws.Range("D11").Value = '107.12'
$ws.Range("E11").Value = '  +0.63%  '

$ws.Range("D12").Value = '1.983.79'
$ws.Range("E12").Value = '  +1.27%  '

$ws.Range("D13").Value = '0.07759'
$ws.Range("E13").Value = '  -0.05%  '

$ws.Range("D14").Value = '5.403'
$ws.Range("E14").Value = '  -0.68%  '

$ws.Range("D15").Value = '0.7062'
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").Value = '279.25'
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").Value = '31.059.13'
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("D18").Value = '13.36'
$ws.Range("E18").Value = '  +1.28%  '

$ws.Range("D19").Value = '0.000007825'
$ws.Range("E19").Value = '  +1.54%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.213.61'
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.536'
$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").Value = '6.522'
$ws.Range("E24").Value = '  -0.90%  '

$ws.Range("D25").Value = '9.783'
$ws.Range("E25").Value = '  -1.77%  '

$ws.Range("D26").Value = '169.34'
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").Value = '19.77'
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("D28").Value = '2.183'
$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").Value = '0.1052'
$ws.Range("E29").Value = '  -0.80%  '

$ws.Range("D30").Value = '1.401'
$ws.Range("E30").Value = '  -2.68%  '

$ws.Range("D31").Value = '4.623'
$ws.Range("E31").Value = '  -3.72%  '

$ws.Range("D32").Value = '1.568'
$ws.Range("E32").Value = '  -1.64%  '

$ws.Range("D33").Value = '4.432'
$ws.Range("E33").Value = '  -1.35%  '

$ws.Range("D34").Value = '0.04914'
$ws.Range("E34").Value = '  -3.15%  '

$ws.Range("D35").Value = '0.7524'
$ws.Range("E35").Value = '  -2.91%  '

$ws.Range("D36").Value = '1.174'
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("D37").Value = '2.733'
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").Value = '0.02012'
$ws.Range("E38").Value = '  -0.98%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.689'
$ws.Range("E39").Value = '  -1.27%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '78.98'
$ws.Range("E40").Value = '  +10.31%  '

$ws.Range("D41").Value = '6.503'
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").Value = '2.137'
$ws.Range("E42").Value = '  +0.91%  '

$ws.Range("D43").Value = '0.8984'
$ws.Range("E43").Value = '  +1.18%  '

$ws.Range("D44").Value = '109.55'
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").Value = '0.4461'
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").Value = '7.923'
$ws.Range("E46").Value = '  +5.59%  '

$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("D48").Value = '988.98'
$ws.Range("E48").Value = '  +4.01%  '

$ws.Range("D49").Value = '0.1251'
$ws.Range("E49").Value = '  -1.43%  '

$ws.Range("D50").Value = '9.337'
$ws.Range("E50").Value = '  -1.32%  '

$ws.Range("D51").Value = '36.05'
$ws.Range("E51").Value = '  +0.47%  '

# Restore the default cell style on column D so no stray number-format
# style survives the edit (matches original workbook formatting).
$priceRange.Style = "Normal"
